$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header "SN" -> "S(N)" in D5
$ws.Range("D5").Value = "S(N)"

# 2. "Zergling x2" row (row 45) becomes just "Zergling" with reduced mineral/supply costs
#    (its Baneling morph target moved into its own new row below)
$ws.Range("A45").Value = "Zergling"
$ws.Range("B45").Value = 25
$ws.Range("F45").Value = 0.5

# 3. Insert a brand-new row for "Baneling" right after the Zergling row,
#    pushing Roach..Ultralisk down by one row.
$ws.Rows.Item(46).Insert()

$ws.Range("A46").Value = "Baneling"
$ws.Range("B46").Value = 50
$ws.Range("C46").Value = 25
$ws.Range("D46").Style = "Normal"
$ws.Range("D46").Value = 44
$ws.Range("E46").Formula = "=D46/1.4"
$ws.Range("F46").Value = 0.5
$ws.Range("G46").Formula = "=((60/E46)*B46)/`$B`$1"
$ws.Range("H46").Formula = "=((60/E46)*(B46+`$F`$3*F46))/`$B`$1"
$ws.Range("I46").Formula = "=((60/E46)*C46)/`$B`$2"

# 4. Extend the two conditional-formatting ranges down to the new last row (57)
$fcs1 = $ws.Range("H6:H56").FormatConditions
$fcs1.Item(1).ModifyAppliesToRange($ws.Range("H6:H57"))

$fcs2 = $ws.Range("I6:I56").FormatConditions
$fcs2.Item(1).ModifyAppliesToRange($ws.Range("I6:I57"))
